# Cadastro de novo funcionário e atualização de disponibilidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Geraldo Magela"
$ws.Range("C2").Value = $true

[void]$ws.Range("A1").Select()
